# Applies the cryptos list price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell reference + new text value.
# "ForceText" = value that looks like a plain number (e.g. "1.00", "0.0490") and
# would otherwise be reinterpreted by Excel as a numeric value (dropping trailing
# zeros / exact decimal representation), so it must be entered as literal text.
$updates = @(
    @{ Cell = "D2"; Value = "55.767.40"; ForceText = $false },
    @{ Cell = "E2"; Value = "  -1.28%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "2.288.96"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -1.16%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  +0.00%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "511.77"; ForceText = $true },
    @{ Cell = "E5"; Value = "  -1.03%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "129.13"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -4.57%  "; ForceText = $false },
    @{ Cell = "D7"; Value = "0.998"; ForceText = $true },
    @{ Cell = "E7"; Value = "  +0.24%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "0.525"; ForceText = $true },
    @{ Cell = "E8"; Value = "  -2.15%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "2.291.03"; ForceText = $false },
    @{ Cell = "E9"; Value = "  -1.89%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.0984"; ForceText = $true },
    @{ Cell = "E10"; Value = "  -3.63%  "; ForceText = $false },
    @{ Cell = "E11"; Value = "  -0.21%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "5.20"; ForceText = $true },
    @{ Cell = "E12"; Value = "  -2.26%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "0.332"; ForceText = $true },
    @{ Cell = "E13"; Value = "  -2.67%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "2.694.65"; ForceText = $false },
    @{ Cell = "E14"; Value = "  -1.29%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "23.02"; ForceText = $true },
    @{ Cell = "E15"; Value = "  -3.80%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "55.781.80"; ForceText = $false },
    @{ Cell = "E16"; Value = "  -1.44%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "0.0000130"; ForceText = $true },
    @{ Cell = "E17"; Value = "  -3.00%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "2.293.40"; ForceText = $false },
    @{ Cell = "E18"; Value = "  -1.76%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "10.26"; ForceText = $true },
    @{ Cell = "E19"; Value = "  -2.11%  "; ForceText = $false },
    @{ Cell = "D20"; Value = "324.95"; ForceText = $true },
    @{ Cell = "E20"; Value = "  +0.65%  "; ForceText = $false },
    @{ Cell = "E21"; Value = "  -3.10%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "6.59"; ForceText = $true },
    @{ Cell = "E22"; Value = "  +0.33%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E23"; Value = "  +0.41%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "60.28"; ForceText = $true },
    @{ Cell = "E24"; Value = "  -0.61%  "; ForceText = $false },
    @{ Cell = "E25"; Value = "  -1.04%  "; ForceText = $false },
    @{ Cell = "B26"; Value = "Binance-PegBSC-USD"; ForceText = $false },
    @{ Cell = "C26"; Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"; ForceText = $false },
    @{ Cell = "D26"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E26"; Value = "  +0.98%  "; ForceText = $false },
    @{ Cell = "B27"; Value = "InternetComputer(DFINITY)"; ForceText = $false },
    @{ Cell = "C27"; Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; ForceText = $false },
    @{ Cell = "D27"; Value = "8.43"; ForceText = $true },
    @{ Cell = "E27"; Value = "  +5.58%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "1.29"; ForceText = $true },
    @{ Cell = "E28"; Value = "  +0.65%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "167.67"; ForceText = $true },
    @{ Cell = "E29"; Value = "  +0.63%  "; ForceText = $false },
    @{ Cell = "E30"; Value = "  -2.07%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "0.0₃0701"; ForceText = $false },
    @{ Cell = "E31"; Value = "  -5.21%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "6.00"; ForceText = $true },
    @{ Cell = "E32"; Value = "  -3.25%  "; ForceText = $false },
    @{ Cell = "B33"; Value = "USDe"; ForceText = $false },
    @{ Cell = "C33"; Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"; ForceText = $false },
    @{ Cell = "D33"; Value = "0.999"; ForceText = $true },
    @{ Cell = "E33"; Value = "  +0.00%  "; ForceText = $false },
    @{ Cell = "B34"; Value = "EthereumClassic"; ForceText = $false },
    @{ Cell = "C34"; Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; ForceText = $false },
    @{ Cell = "D34"; Value = "18.08"; ForceText = $true },
    @{ Cell = "E34"; Value = "  -1.59%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "0.999"; ForceText = $true },
    @{ Cell = "E35"; Value = "  +0.70%  "; ForceText = $false },
    @{ Cell = "E36"; Value = "  -3.09%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "3.86"; ForceText = $true },
    @{ Cell = "E37"; Value = "  -3.84%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "0.875"; ForceText = $true },
    @{ Cell = "E38"; Value = "  -5.16%  "; ForceText = $false },
    @{ Cell = "B39"; Value = "OKB"; ForceText = $false },
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; ForceText = $false },
    @{ Cell = "D39"; Value = "38.34"; ForceText = $true },
    @{ Cell = "E39"; Value = "  +1.06%  "; ForceText = $false },
    @{ Cell = "B40"; Value = "Stacks"; ForceText = $false },
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; ForceText = $false },
    @{ Cell = "D40"; Value = "1.54"; ForceText = $true },
    @{ Cell = "E40"; Value = "  -0.25%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "145.65"; ForceText = $true },
    @{ Cell = "E41"; Value = "  +3.98%  "; ForceText = $false },
    @{ Cell = "E42"; Value = "  -3.34%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "3.52"; ForceText = $true },
    @{ Cell = "E43"; Value = "  -2.43%  "; ForceText = $false },
    @{ Cell = "D44"; Value = "279.10"; ForceText = $true },
    @{ Cell = "E44"; Value = "  +0.99%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "4.86"; ForceText = $true },
    @{ Cell = "E45"; Value = "  -7.50%  "; ForceText = $false },
    @{ Cell = "D46"; Value = "0.0918"; ForceText = $true },
    @{ Cell = "E46"; Value = "  -1.55%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "0.0490"; ForceText = $true },
    @{ Cell = "E47"; Value = "  -3.39%  "; ForceText = $false },
    @{ Cell = "D48"; Value = "0.548"; ForceText = $true },
    @{ Cell = "E48"; Value = "  -2.11%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "17.81"; ForceText = $true },
    @{ Cell = "E49"; Value = "  +0.24%  "; ForceText = $false },
    @{ Cell = "B50"; Value = "Polygon"; ForceText = $false },
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; ForceText = $false },
    @{ Cell = "D50"; Value = "0.377"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -1.08%  "; ForceText = $false },
    @{ Cell = "B51"; Value = "VeChain"; ForceText = $false },
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; ForceText = $false },
    @{ Cell = "D51"; Value = "0.0211"; ForceText = $true },
    @{ Cell = "E51"; Value = "  -2.84%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Leading apostrophe forces Excel to store the value as text instead of a number.
        $range.Value = "'" + $u.Value
        # Re-apply the default "Normal" style so the text-entry does not leave behind
        # a Text number-format style that the cell did not have before.
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
